$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying the existing "2022-Q3" sheet
#    (so it inherits the same column styles / header formatting), placed
#    immediately before it. This also shifts 2022-Q3 / 2021-Q4 / 2021-Q3
#    one position to the right, matching the target sheet order:
#      总计, 2022-Q4, 2022-Q3, 2021-Q4, 2021-Q3
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Copy($q3Sheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Wipe the copied values (keep formatting) before writing fresh data.
$q4Sheet.Range("A1:H20").ClearContents()

# Propagate the index-column / header style (taken from the still-styled
# A2 cell after ClearContents) down across all the rows we are about to fill.
$q4Sheet.Range("A2").Copy()
$q4Sheet.Range("A2:A8").PasteSpecial(-4122)

# Headers
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Force text formatting on the fund-code / text-looking numeric columns so
# leading zeros and exact decimal text are preserved (matches source: these
# columns are stored as text, not numbers).
$q4Sheet.Range("B2:B8").NumberFormat = "@"
$q4Sheet.Range("D2:G8").NumberFormat = "@"

$q4Data = @(
    @(0, "009686", "华夏磐利一年定期开放混合A",   "11.49", "92.56", "6.00", "0.6894", 1),
    @(1, "015697", "华夏磐润两年定开混合A",       "2.68",  "86.76", "4.68", "0.1254", 3),
    @(2, "015698", "华夏磐润两年定开混合C",       "0.99",  "86.76", "4.68", "0.0463", 3),
    @(3, "001339", "兴银鼎新灵活配置混合",         "0.71",  "87.16", "5.14", "0.0365", 3),
    @(4, "009687", "华夏磐利一年定期开放混合C",   "0.46",  "92.56", "6.00", "0.0276", 1),
    @(5, "010124", "兴银景气优选混合A",           "0.45",  "83.23", "5.10", "0.0230", 1),
    @(6, "010125", "兴银景气优选混合C",           "0.35",  "83.23", "5.10", "0.0178", 1)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = 2 + $i
    $rec = $q4Data[$i]
    $q4Sheet.Cells.Item($row, 1).Value = $rec[0]
    $q4Sheet.Cells.Item($row, 2).Value = $rec[1]
    $q4Sheet.Cells.Item($row, 3).Value = $rec[2]
    $q4Sheet.Cells.Item($row, 4).Value = $rec[3]
    $q4Sheet.Cells.Item($row, 5).Value = $rec[4]
    $q4Sheet.Cells.Item($row, 6).Value = $rec[5]
    $q4Sheet.Cells.Item($row, 7).Value = $rec[6]
    $q4Sheet.Cells.Item($row, 8).Value = $rec[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the new 2022-Q4 row at the top
#    of the data (row 2) and push the previously existing rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Row 5 is brand new (the old data only went to row 4) -- extend the
# index-column style down to it so it matches its siblings (s="2").
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$totalData = @(
    @(0, "2022-Q4", 7, 0.97),
    @(1, "2022-Q3", 2, 0.02),
    @(2, "2021-Q4", 6, 1.68),
    @(3, "2021-Q3", 1, 1.62)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $row = 2 + $i
    $rec = $totalData[$i]
    $total.Cells.Item($row, 1).Value = $rec[0]
    $total.Cells.Item($row, 2).Value = $rec[1]
    $total.Cells.Item($row, 3).Value = $rec[2]
    $total.Cells.Item($row, 4).Value = $rec[3]
}

# Keep the original active tab (总计 / sheet index 1) selected, matching the
# workbook's unchanged bookViews activeTab="0".
$total.Activate()
